# Update "Des Scheduled Flights vs actual.xlsx"
# - correct the 2021-10-25 row's actual-flights count (C568: 70 -> 71)
# - append 20 new daily rows (2021-10-26 .. 2021-11-14) with scheduled /
#   actual flight counts and the C/B completion-rate formula
# - leave the selection on B600 (mirrors the author's final cursor spot)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- fix existing row 568 (2021-10-25): actual flights 70 -> 71 ----
$ws.Range("C568").Value = 71

# ---- new data rows 569..588 ----
# columns: row, date text, B scheduled, C actual
$rows = @(
    @(569, "2021-10-26", 73, 68),
    @(570, "2021-10-27", 64, 62),
    @(571, "2021-10-28", 69, 65),
    @(572, "2021-10-29", 82, 76),
    @(573, "2021-10-30", 59, 57),
    @(574, "2021-10-31", 48, 47),
    @(575, "2021-11-01", 69, 63),
    @(576, "2021-11-02", 68, 64),
    @(577, "2021-11-03", 61, 54),
    @(578, "2021-11-04", 78, 73),
    @(579, "2021-11-05", 75, 72),
    @(580, "2021-11-06", 59, 57),
    @(581, "2021-11-07", 61, 57),
    @(582, "2021-11-08", 67, 65),
    @(583, "2021-11-09", 65, 64),
    @(584, "2021-11-10", 71, 66),
    @(585, "2021-11-11", 84, 78),
    @(586, "2021-11-12", 81, 72),
    @(587, "2021-11-13", 63, 60),
    @(588, "2021-11-14", 60, 58)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $dateText = $r[1]
    $scheduled = $r[2]
    $actual = $r[3]

    # Copy the formatting of the row directly above so the new row
    # keeps the same styles (text-centered date column, number columns,
    # percentage formula column) as the rest of the table.
    $ws.Range("A$($rowNum - 1):D$($rowNum - 1)").Copy()
    $ws.Range("A$($rowNum):D$($rowNum)").PasteSpecial(-4122)

    $ws.Range("A$rowNum").NumberFormat = "@"
    $ws.Range("A$rowNum").Value = $dateText
    $ws.Range("B$rowNum").Value = $scheduled
    $ws.Range("C$rowNum").Value = $actual
    $ws.Range("D$rowNum").Formula = "=C$rowNum/B$rowNum"
}

$excel.CutCopyMode = $false

# ---- restore the view to where the author left it ----
$ws.Range("B600").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 347
$win.ScrollColumn = 1
